{"js": "// Office.js (Word JavaScript API) script.\n// Adds the \"Historias de Usuarios\" content: a centered bold title followed\n// by six user-story entries (T\u00edtulo / Descripci\u00f3n / Condiciones de\n// aceptaci\u00f3n), inserted before the document's original (sole, empty)\n// paragraph. The insertOoxml package below contains 25 paragraphs; Word\n// merges the content of the *last* one into the pre-existing paragraph so\n// that its bookmark (\"_GoBack\") ends up on the final \"Condiciones de\n// aceptaci\u00f3n\" paragraph, matching the target document exactly.\n\nconst body = context.document.body;\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p><w:pPr><w:jc w:val=\"center\"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>HISTORIAS DE USUARIOS</w:t></w:r></w:p><w:p><w:r><w:t>1)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>T\u00edtulo:</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Alta de una reserva</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Descripci\u00f3n: </w:t></w:r><w:r><w:t>Q</w:t></w:r><w:r><w:t>uisiera poder cargar una reserva, siempre y cuando las condiciones est\u00e9n dadas, es decir haya salones con mesas disponibles, mozos disponibles, etc\u2026</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condiciones de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p><w:p><w:r><w:t>2)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">T\u00edtulo: </w:t></w:r><w:r><w:t>Alta de un mozo</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Descripci\u00f3n: </w:t></w:r><w:r><w:t>Q</w:t></w:r><w:r><w:t>uisiera poder cargar un mozo siempre y cuando no haya 6 mozos en el restaurante.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condici\u00f3n de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p><w:p><w:r><w:t>3)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">T\u00edtulo: </w:t></w:r><w:r><w:t>Finalizar Reserva</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Descripci\u00f3n</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r><w:r><w:t>Quisiera poder finalizar una reserva ingresando su id y luego poder ingresar el monto que le corresponde a pagar al cliente</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condiciones de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p><w:p><w:r><w:t>4)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">T\u00edtulo: </w:t></w:r><w:r><w:t>Mostrar disponibilidad de un sal\u00f3n</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Descripci\u00f3n: </w:t></w:r><w:r><w:t>Quisiera poder mostrar la cantidad de mesas que le quedan disponibles a un sal\u00f3n ingresando su nombre</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condiciones de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p><w:p><w:r><w:t>5)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">T\u00edtulo: </w:t></w:r><w:r><w:t>Listar reservas</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Descripci\u00f3n: </w:t></w:r><w:r><w:t>Quisiera poder listar las reservas del restaurantes mostrando su correspondiente informaci\u00f3n (en caso del sal\u00f3n, el cliente y el mozo solo mostrar su nombre)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condiciones de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p><w:p><w:r><w:t>6)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">T\u00edtulo: </w:t></w:r><w:r><w:t>Eliminar una reserva ingresando su id</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Descripci\u00f3n: </w:t></w:r><w:r><w:t>Quisiera poder eliminar (de la base de datos) una reserva ingresando su id</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condiciones de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nbody.insertOoxml(ooxml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Adds the \"Historias de Usuarios\" content: a centered bold title followed\n# by six user-story entries (Titulo / Descripcion / Condiciones de\n# aceptacion), inserted before the document's original (sole, empty)\n# paragraph. $xml below packages 25 paragraphs; Range.InsertXML merges the\n# content of the *last* one into the collapsed range's paragraph, so the\n# pre-existing paragraph's bookmark (\"_GoBack\") ends up on the final\n# \"Condiciones de aceptacion\" paragraph - matching the target document.\n\n$d = $word.ActiveDocument\n\n$r = $d.Paragraphs(1).Range\n$r.Collapse(1)  # wdCollapseStart\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p><w:pPr><w:jc w:val=\"center\"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>HISTORIAS DE USUARIOS</w:t></w:r></w:p><w:p><w:r><w:t>1)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>T\u00edtulo:</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Alta de una reserva</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Descripci\u00f3n: </w:t></w:r><w:r><w:t>Q</w:t></w:r><w:r><w:t>uisiera poder cargar una reserva, siempre y cuando las condiciones est\u00e9n dadas, es decir haya salones con mesas disponibles, mozos disponibles, etc\u2026</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condiciones de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p><w:p><w:r><w:t>2)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">T\u00edtulo: </w:t></w:r><w:r><w:t>Alta de un mozo</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Descripci\u00f3n: </w:t></w:r><w:r><w:t>Q</w:t></w:r><w:r><w:t>uisiera poder cargar un mozo siempre y cuando no haya 6 mozos en el restaurante.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condici\u00f3n de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p><w:p><w:r><w:t>3)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">T\u00edtulo: </w:t></w:r><w:r><w:t>Finalizar Reserva</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Descripci\u00f3n</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r><w:r><w:t>Quisiera poder finalizar una reserva ingresando su id y luego poder ingresar el monto que le corresponde a pagar al cliente</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condiciones de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p><w:p><w:r><w:t>4)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">T\u00edtulo: </w:t></w:r><w:r><w:t>Mostrar disponibilidad de un sal\u00f3n</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Descripci\u00f3n: </w:t></w:r><w:r><w:t>Quisiera poder mostrar la cantidad de mesas que le quedan disponibles a un sal\u00f3n ingresando su nombre</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condiciones de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p><w:p><w:r><w:t>5)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">T\u00edtulo: </w:t></w:r><w:r><w:t>Listar reservas</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Descripci\u00f3n: </w:t></w:r><w:r><w:t>Quisiera poder listar las reservas del restaurantes mostrando su correspondiente informaci\u00f3n (en caso del sal\u00f3n, el cliente y el mozo solo mostrar su nombre)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condiciones de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p><w:p><w:r><w:t>6)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">T\u00edtulo: </w:t></w:r><w:r><w:t>Eliminar una reserva ingresando su id</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Descripci\u00f3n: </w:t></w:r><w:r><w:t>Quisiera poder eliminar (de la base de datos) una reserva ingresando su id</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">Condiciones de aceptaci\u00f3n: </w:t></w:r><w:r><w:t>Lo que se solicita</w:t></w:r></w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>'\n\n$r.InsertXML($xml)\n"}
